$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H33").Value = 97.69231000000001
$ws.Range("I33").Value = 79.09090999999999
$ws.Range("K33").Value = 79.09090999999999
$ws.Range("M33").Value = 149.90909
$ws.Range("H46").Value = 6316
$ws.Range("J46").Value = 6316
$ws.Range("L46").Value = 18948
$ws.Range("N46").Value = -19186
$ws.Range("H60").Value = 6316
$ws.Range("J60").Value = 6316
$ws.Range("L60").Value = 18948
$ws.Range("N60").Value = -19916
$ws.Range("H62").Value = 5905.8887
$ws.Range("J62").Value = 9362.375
$ws.Range("L62").Value = 9362.375
$ws.Range("N62").Value = -10610.375
$ws.Range("H64").Value = 5222.6665
$ws.Range("I64").Value = 3940
$ws.Range("J64").Value = 8750
$ws.Range("K64").Value = 3940
$ws.Range("L64").Value = 8750
$ws.Range("M64").Value = -3692
$ws.Range("N64").Value = -9246
$ws.Range("H65").Value = 5905.8887
$ws.Range("J65").Value = 9362.375
$ws.Range("L65").Value = 46811.875
$ws.Range("N65").Value = -53051.875
$ws.Range("H67").Value = 5222.6665
$ws.Range("I67").Value = 3940
$ws.Range("J67").Value = 8750
$ws.Range("K67").Value = 3940
$ws.Range("L67").Value = 8750
$ws.Range("M67").Value = -3082
$ws.Range("N67").Value = -10466
$ws.Range("H74").Value = 6299.4443
$ws.Range("I74").Value = 5528
$ws.Range("K74").Value = 5528
$ws.Range("M74").Value = -4592
$ws.Range("H76").Value = 0
$ws.Range("I76").Value = 0
$ws.Range("J76").Value = 0
$ws.Range("K76").Value = 0
$ws.Range("L76").ClearContents()
$ws.Range("M76").ClearContents()
$ws.Range("N76").Value = 0
$ws.Range("H77").Value = 6299.4443
$ws.Range("I77").Value = 5528
$ws.Range("K77").Value = 27640
$ws.Range("M77").Value = -22960
$ws.Range("H79").Value = 0
$ws.Range("I79").Value = 0
$ws.Range("J79").Value = 0
$ws.Range("K79").Value = 0
$ws.Range("L79").ClearContents()
$ws.Range("M79").ClearContents()
$ws.Range("N79").Value = 0
$ws.Range("H86").Value = 0
$ws.Range("I86").Value = 0
$ws.Range("K86").Value = 0
$ws.Range("M86").ClearContents()
$ws.Range("H89").Value = 0
$ws.Range("I89").Value = 0
$ws.Range("K89").Value = 0
$ws.Range("M89").ClearContents()
$ws.Range("H113").Value = 4999.5
$ws.Range("I113").Value = 4999.5
$ws.Range("J113").Value = 0
$ws.Range("K113").Value = 4999.5
$ws.Range("L113").Value = 0
$ws.Range("M113").ClearContents()
$ws.Range("N113").Value = -1745.5
$ws.Range("H137").Value = 1583.2858
$ws.Range("I137").Value = 1148.75
$ws.Range("K137").Value = 3446.25
$ws.Range("M137").Value = -896.25
$ws.Range("H138").Value = 3265.0557
$ws.Range("J138").Value = 3988.5454
$ws.Range("L138").Value = 11965.6362
$ws.Range("N138").Value = -22245.6362

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H4").Value = 599.4
$ws.Range("I4").Value = 599.4
$ws.Range("K4").Value = 599.4
$ws.Range("M4").Value = -483.4
$ws.Range("H88").Value = 603.2222
$ws.Range("J88").Value = 782.8
$ws.Range("L88").Value = 782.8
$ws.Range("N88").Value = -1594.8
$ws.Range("H91").Value = 603.2222
$ws.Range("J91").Value = 782.8
$ws.Range("L91").Value = 782.8
$ws.Range("N91").Value = -3590.8
$ws.Range("H135").Value = 81689.25
$ws.Range("J135").Value = 81689.25
$ws.Range("L135").Value = 81689.25
$ws.Range("N135").Value = -91829.25

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H9").Value = 20000
$ws.Range("J9").Value = 20000
$ws.Range("L9").Value = 20000
$ws.Range("N9").Value = -20336
$ws.Range("H64").Value = 733.5714
$ws.Range("I64").Value = 713.25
$ws.Range("J64").Value = 760.6667
$ws.Range("K64").Value = 713.25
$ws.Range("L64").Value = 760.6667
$ws.Range("M64").Value = -488.25
$ws.Range("N64").Value = -1210.6667
$ws.Range("H67").Value = 733.5714
$ws.Range("I67").Value = 713.25
$ws.Range("J67").Value = 760.6667
$ws.Range("K67").Value = 713.25
$ws.Range("L67").Value = 760.6667
$ws.Range("M67").Value = 66.75
$ws.Range("N67").Value = -2320.6667
$ws.Range("H86").Value = 4140.636
$ws.Range("I86").Value = 1756.6666
$ws.Range("K86").Value = 1756.6666
$ws.Range("M86").Value = -633.6666
$ws.Range("H89").Value = 4140.636
$ws.Range("I89").Value = 1756.6666
$ws.Range("K89").Value = 8783.333000000001
$ws.Range("M89").Value = -3167.333000000001
$ws.Range("H94").Value = 449.66666
$ws.Range("I94").Value = 349.625
$ws.Range("K94").Value = 349.625
$ws.Range("M94").Value = 101.375
$ws.Range("H95").Value = 5764
$ws.Range("J95").Value = 5764
$ws.Range("L95").Value = 5764
$ws.Range("N95").Value = -11256
$ws.Range("H105").Value = 2055
$ws.Range("I105").Value = 2055
$ws.Range("K105").Value = 2055
$ws.Range("M105").Value = -308
$ws.Range("H134").Value = 0
$ws.Range("I134").Value = 0
$ws.Range("K134").Value = 0
$ws.Range("M134").ClearContents()

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 6972.769
$ws.Range("J31").Value = 8428.056
$ws.Range("L31").Value = 8428.056
$ws.Range("N31").Value = -9018.056
$ws.Range("H34").Value = 6972.769
$ws.Range("J34").Value = 8428.056
$ws.Range("L34").Value = 8428.056
$ws.Range("N34").Value = -8832.056
$ws.Range("H59").Value = 40441.25
$ws.Range("I59").Value = 15000
$ws.Range("J59").Value = 48921.668
$ws.Range("K59").Value = 15000
$ws.Range("L59").Value = 48921.668
$ws.Range("M59").Value = -13855
$ws.Range("N59").Value = -51211.668
$ws.Range("H134").Value = 5000
$ws.Range("I134").Value = 5000
$ws.Range("K134").Value = 15000
$ws.Range("M134").Value = -12465

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H131").Value = 2661.9443
$ws.Range("I131").Value = 2306
$ws.Range("K131").Value = 6918
$ws.Range("M131").Value = -1878
$ws.Range("H139").Value = 3720
$ws.Range("I139").Value = 3327.2222
$ws.Range("K139").Value = 9981.6666
$ws.Range("M139").Value = -4841.6666
$ws.Range("H140").Value = 3922.9167
$ws.Range("I140").Value = 3657.5
$ws.Range("J140").Value = 5250
$ws.Range("K140").Value = 10972.5
$ws.Range("L140").Value = 15750
$ws.Range("M140").Value = -5792.5
$ws.Range("N140").Value = -26110

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H97").Value = 935.44446
$ws.Range("I97").Value = 727.75
$ws.Range("K97").Value = 727.75
$ws.Range("M97").Value = -231.75
$ws.Range("H113").Value = 5831.75
$ws.Range("I113").Value = 3529.8
$ws.Range("K113").Value = 3529.8
$ws.Range("M113").Value = -1359.8
$ws.Range("H132").Value = 75501
$ws.Range("I132").Value = 103371.4
$ws.Range("J132").Value = 5825
$ws.Range("K132").Value = 310114.2
$ws.Range("L132").Value = 17475
$ws.Range("M132").Value = -307584.2
$ws.Range("N132").Value = -22535

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H55").Value = 735.61536
$ws.Range("I55").Value = 1440
$ws.Range("J55").Value = 295.375
$ws.Range("K55").Value = 1440
$ws.Range("L55").Value = 295.375
$ws.Range("M55").Value = -1267
$ws.Range("N55").Value = -641.375
$ws.Range("H82").Value = 3153.5
$ws.Range("J82").Value = 4525.7144
$ws.Range("L82").Value = 4525.7144
$ws.Range("N82").Value = -5247.7144
$ws.Range("H85").Value = 3153.5
$ws.Range("J85").Value = 4525.7144
$ws.Range("L85").Value = 4525.7144
$ws.Range("N85").Value = -7021.7144
$ws.Range("H93").Value = 1599.8
$ws.Range("I93").Value = 1550.375
$ws.Range("K93").Value = 1550.375
$ws.Range("M93").Value = -302.375

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H64").Value = 0
$ws.Range("J64").Value = 0
$ws.Range("L64").ClearContents()
$ws.Range("N64").Value = 0
$ws.Range("H67").Value = 0
$ws.Range("J67").Value = 0
$ws.Range("L67").ClearContents()
$ws.Range("N67").Value = 0
$ws.Range("H107").Value = 835.9091
$ws.Range("I107").Value = 516
$ws.Range("J107").Value = 1219.8
$ws.Range("K107").Value = 1548
$ws.Range("L107").Value = 3659.4
$ws.Range("M107").Value = 372
$ws.Range("N107").Value = -7499.4
$ws.Range("H132").Value = 1674.6666
$ws.Range("I132").Value = 1265
$ws.Range("J132").Value = 2801.25
$ws.Range("K132").Value = 3795
$ws.Range("L132").Value = 8403.75
$ws.Range("M132").Value = -1265
$ws.Range("N132").Value = -13463.75
